$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8189005255699158
$ws.Range("B1").Value = 3.047157049179077
$ws.Range("C1").Value = 2.597966909408569
$ws.Range("D1").Value = 2.262613296508789
$ws.Range("E1").Value = 1.916884779930115
